# Auto-generated script to apply cryptos.xlsx price/volume update
# Commit: "Updated cryptos list on Tue Oct 31 08:53:41 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.436.21"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.801.79"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "227.60"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "0.574"
$ws.Range("E6").Value = "  +3.13%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "36.21"
$ws.Range("E8").Value = "  +7.25%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "0.0963"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "2.059.64"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "11.58"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.797.31"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.645"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "4.50"
$ws.Range("E16").Value = "  +4.85%  "
$ws.Range("D17").Value = "34.391.32"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "68.99"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "244.81"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").Value = "11.60"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "172.54"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").Value = "7.97"
$ws.Range("E26").Value = "  +8.57%  "
$ws.Range("D27").Value = "16.83"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "0.118"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").Value = "0.0528"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "3.84"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "1.397.36"
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("D36").Value = "0.672"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  -6.44%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "82.47"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").Value = "0.959"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("E44").Value = "  +6.39%  "
$ws.Range("D45").Value = "13.33"
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("D46").Value = "6.03"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "0.0502"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").Value = "1.962.27"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "104.37"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "0.0₆0123"
$ws.Range("E51").Value = "  -1.86%  "
